$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (Beta) values ---
$ws.Cells.Item(2,3).Value  = 19.12075701903682   # C2
$ws.Cells.Item(2,5).Value  = 0.01982943797740053 # E2
$ws.Cells.Item(2,6).Value  = 10.34362016733451   # F2
$ws.Cells.Item(2,7).Value  = 9.90645192253899    # G2
$ws.Cells.Item(2,8).Value  = 10.7966804998958    # H2
$ws.Cells.Item(2,9).Value  = 0.003147955924683531 # I2
$ws.Cells.Item(2,10).Value = 0.002757539016402122 # J2
$ws.Cells.Item(2,11).Value = 0.003631286706931039 # K2
$ws.Cells.Item(2,12).Value = 0.01105691775784921  # L2
$ws.Cells.Item(2,13).Value = 0.01073834377938824  # M2
$ws.Cells.Item(2,14).Value = 0.01139049755384476  # N2

# --- Update existing row 3 (Gamma) values ---
$ws.Cells.Item(3,3).Value  = 0.04981522627320694  # C3
$ws.Cells.Item(3,4).Value  = 0.04815098319456564  # D3
$ws.Cells.Item(3,5).Value  = 0.0499839736740351   # E3
$ws.Cells.Item(3,6).Value  = 0.1137245776403948   # F3
$ws.Cells.Item(3,7).Value  = 0.03237751211294077  # G3
$ws.Cells.Item(3,8).Value  = 0.2124944064076332   # H3
$ws.Cells.Item(3,9).Value  = 0.1051142188779256   # I3
$ws.Cells.Item(3,10).Value = 0.03076849881608078  # J3
$ws.Cells.Item(3,11).Value = 0.1954132764324292   # K3
$ws.Cells.Item(3,12).Value = 0.1201171135361337   # L3
$ws.Cells.Item(3,13).Value = 0.03315277237719329  # M3
$ws.Cells.Item(3,14).Value = 0.2256327657098036   # N3

# --- Add new row 4 (Beta + Gamma) ---
# Copy the style from A3 (bold/border/center) onto A4 before setting values
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)

$ws.Cells.Item(4,1).Value  = 2
$ws.Cells.Item(4,2).Value  = "Beta + Gamma"
$ws.Cells.Item(4,3).Value  = 19.17057224531003
$ws.Cells.Item(4,4).Value  = 0.05509537454402212
$ws.Cells.Item(4,5).Value  = 0.06981341165143562
$ws.Cells.Item(4,6).Value  = 10.4573447449749
$ws.Cells.Item(4,7).Value  = 9.938829434651929
$ws.Cells.Item(4,8).Value  = 11.00917490630343
$ws.Cells.Item(4,9).Value  = 0.1082621748026091
$ws.Cells.Item(4,10).Value = 0.0335260378324829
$ws.Cells.Item(4,11).Value = 0.1990445631393603
$ws.Cells.Item(4,12).Value = 0.1311740312939829
$ws.Cells.Item(4,13).Value = 0.04389111615658152
$ws.Cells.Item(4,14).Value = 0.2370232632636484
